$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (46081 -> 46082) for every data row (rows 2 through 128).
$ws.Range("C2:C128").Value = 46082
